$wb = $excel.ActiveWorkbook

# --- 1. Add a "State" column to the hotel_info sheet, between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Range("C1").EntireColumn.Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheets: review_info first, hotel_info second ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
